$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1524.4117
$ws.Range("J17").Value = 1524.4117
$ws.Range("L17").Value = 4573.2351
$ws.Range("N17").Value = -4909.2351
$ws.Range("H82").Value = 2750
$ws.Range("I82").Value = 2750
$ws.Range("K82").Value = 8250
$ws.Range("M82").Value = -7844
$ws.Range("H85").Value = 2750
$ws.Range("I85").Value = 2750
$ws.Range("K85").Value = 8250
$ws.Range("M85").Value = -6846
$ws.Range("H87").Value = 44996.25
$ws.Range("J87").Value = 44996.25
$ws.Range("L87").Value = 44996.25
$ws.Range("N87").Value = -47492.25
$ws.Range("H88").Value = 13601.125
$ws.Range("I88").Value = 2750
$ws.Range("J88").Value = 17218.166
$ws.Range("K88").Value = 2750
$ws.Range("L88").Value = 17218.166
$ws.Range("M88").Value = -2344
$ws.Range("N88").Value = -18030.166
$ws.Range("H90").Value = 44996.25
$ws.Range("J90").Value = 44996.25
$ws.Range("L90").Value = 134988.75
$ws.Range("N90").Value = -147468.75
$ws.Range("H91").Value = 13601.125
$ws.Range("I91").Value = 2750
$ws.Range("J91").Value = 17218.166
$ws.Range("K91").Value = 2750
$ws.Range("L91").Value = 17218.166
$ws.Range("M91").Value = -1346
$ws.Range("N91").Value = -20026.166
$ws.Range("H100").Value = 3059.4443
$ws.Range("I100").Value = 3279.4
$ws.Range("J100").Value = 2784.5
$ws.Range("K100").Value = 3279.4
$ws.Range("L100").Value = 2784.5
$ws.Range("M100").Value = -2738.4
$ws.Range("N100").Value = -3866.5
$ws.Range("H111").Value = 488.81818
$ws.Range("I111").Value = 451.66666
$ws.Range("K111").Value = 1354.99998
$ws.Range("M111").Value = 1712.00002
$ws.Range("H129").Value = 612
$ws.Range("I129").Value = 612
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1836
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 3164
$ws.Range("N129").ClearContents()
$ws.Range("H137").Value = 3742.6
$ws.Range("I137").Value = 3645.5
$ws.Range("J137").Value = 3888.25
$ws.Range("K137").Value = 10936.5
$ws.Range("L137").Value = 11664.75
$ws.Range("M137").Value = -8386.5
$ws.Range("N137").Value = -16764.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 56224.5
$ws.Range("J2").Value = 109949
$ws.Range("L2").Value = 109949
$ws.Range("N2").Value = -110175
$ws.Range("H74").Value = 3378.25
$ws.Range("I74").Value = 3066.3333
$ws.Range("K74").Value = 3066.3333
$ws.Range("M74").Value = -2192.3333
$ws.Range("H77").Value = 3378.25
$ws.Range("I77").Value = 3066.3333
$ws.Range("K77").Value = 15331.6665
$ws.Range("M77").Value = -10963.6665
$ws.Range("H116").Value = 56224.5
$ws.Range("J116").Value = 109949
$ws.Range("L116").Value = 109949
$ws.Range("N116").Value = -114537
$ws.Range("H132").Value = 2788.7778
$ws.Range("I132").Value = 2788.7778
$ws.Range("K132").Value = 8366.3334
$ws.Range("M132").Value = -5836.3334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 56224.5
$ws.Range("J3").Value = 109949
$ws.Range("L3").Value = 109949
$ws.Range("N3").Value = -110177
$ws.Range("H107").Value = 780.4286
$ws.Range("I107").Value = 780.4286
$ws.Range("K107").Value = 780.4286
$ws.Range("M107").Value = 1139.5714
$ws.Range("H138").Value = 49999
$ws.Range("J138").Value = 49999
$ws.Range("L138").Value = 49999
$ws.Range("N138").Value = -60279

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2500
$ws.Range("J16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("N16").Value = -3074
$ws.Range("H31").Value = 1705.9166
$ws.Range("I31").Value = 1705.9166
$ws.Range("K31").Value = 1705.9166
$ws.Range("M31").Value = -1410.9166
$ws.Range("H34").Value = 1705.9166
$ws.Range("I34").Value = 1705.9166
$ws.Range("K34").Value = 1705.9166
$ws.Range("M34").Value = -1503.9166
$ws.Range("H50").Value = 20092.092
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -21250
$ws.Range("H59").Value = 25708.584
$ws.Range("J59").Value = 31249.875
$ws.Range("L59").Value = 31249.875
$ws.Range("N59").Value = -33539.875
$ws.Range("H107").Value = 737.7368
$ws.Range("J107").Value = 678.5
$ws.Range("L107").Value = 678.5
$ws.Range("N107").Value = -4518.5
$ws.Range("H113").Value = 2500
$ws.Range("J113").Value = 2500
$ws.Range("L113").Value = 2500
$ws.Range("N113").Value = -6840
$ws.Range("H122").Value = 2027.625
$ws.Range("I122").Value = 1960.1428
$ws.Range("K122").Value = 5880.428400000001
$ws.Range("M122").Value = -3430.428400000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 100000
$ws.Range("J37").Value = 100000
$ws.Range("L37").Value = 300000
$ws.Range("N37").Value = -300224

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1596.2222
$ws.Range("I97").Value = 545.75
$ws.Range("K97").Value = 545.75
$ws.Range("M97").Value = -49.75
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 1699
$ws.Range("I122").Value = 1699
$ws.Range("K122").Value = 5097
$ws.Range("M122").Value = -2647
$ws.Range("H132").Value = 2499.5
$ws.Range("I132").Value = 2499.5
$ws.Range("K132").Value = 7498.5
$ws.Range("M132").Value = -4968.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4190
$ws.Range("I40").Value = 4373.846
$ws.Range("J40").Value = 1800
$ws.Range("K40").Value = 4373.846
$ws.Range("L40").Value = 1800
$ws.Range("M40").Value = -4237.846
$ws.Range("N40").Value = -2072
$ws.Range("H46").Value = 3480.4614
$ws.Range("I46").Value = 2963.7144
$ws.Range("J46").Value = 4083.3333
$ws.Range("K46").Value = 2963.7144
$ws.Range("L46").Value = 4083.3333
$ws.Range("M46").Value = -2775.7144
$ws.Range("N46").Value = -4459.3333
$ws.Range("H61").Value = 1649.909
$ws.Range("I61").Value = 1759.7778
$ws.Range("K61").Value = 1759.7778
$ws.Range("M61").Value = -1557.7778
$ws.Range("H113").Value = 1649.909
$ws.Range("I113").Value = 1759.7778
$ws.Range("K113").Value = 1759.7778
$ws.Range("M113").Value = 410.2221999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 616000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 616000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 616000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -617040
$ws.Range("H81").Value = 1790
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1790
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H97").Value = 67786
$ws.Range("J97").Value = 67786
$ws.Range("L97").Value = 67786
$ws.Range("N97").Value = -69768
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990
$ws.Range("H132").Value = 5000
$ws.Range("J132").Value = 5000
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -20060

Write-Output "Applied all Marilith_Profits updates"